# add test case Car_Reservation_04
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Car_Reservation_04"
$ws.Range("B9").Value = "Car_SRS_27"
$ws.Range("C9").Value = "Functional"
$ws.Range("D9").Value = "Check if the date of reservation is correct after reserve a car"
$ws.Range("E9").Value = "1) open URL ""http://CarPurchasing""`n2)Login with an existing account"
$ws.Range("F9").Value = "User name: customer`npassword: soso.soso126"
$ws.Range("H9").Value = "User should be redirected to Reserved cars page`nand this page contains all the cars reserved by the user before which written on it reserved by username"
$ws.Range("G9").Value = "1)From the home page click on ""see more"" button at any car`n2)Click on Reserve button `n3)From the header click on ""Reserved cars"" `n4)check if the car exists and the data is right on it "
$ws.Range("J9").Value = "Fatma"
$ws.Range("K9").Value = "passed"
$ws.Range("M9").Value = "jannat"

$ws.Rows.Item(9).RowHeight = 131.25

$ws.Range("G9").Select()
